# Fix maintenance log detail views, add better layout, and translations
$wb = $excel.ActiveWorkbook

# --- "initial" sheet: move the selection, it is no longer the active tab ---
$wsInitial = $wb.Worksheets.Item("initial")
$wsInitial.Range("P14").Select()

# --- "common_translations" sheet: add the new spare-parts translation rows ---
$wsCommon = $wb.Worksheets.Item("common_translations")

# Columns A & B (token + English text) are filled in first, row by row ...
$wsCommon.Range("A151").Value = "spare_parts_electrical"
$wsCommon.Range("B151").Value = "Spare Parts (Electrical):"
$wsCommon.Range("A152").Value = "spare_parts_hardware"
$wsCommon.Range("B152").Value = "Spare Parts (Hardware):"
$wsCommon.Range("A153").Value = "spare_parts_monitoring"
$wsCommon.Range("B153").Value = "Spare Parts (Monitoring):"
$wsCommon.Range("A154").Value = "spare_parts_power"
$wsCommon.Range("B154").Value = "Spare Parts (Power):"
$wsCommon.Range("A155").Value = "spare_parts_refrigeration"
$wsCommon.Range("B155").Value = "Spare Parts (Refrigeration):"
$wsCommon.Range("A156").Value = "spare_parts_solar"
$wsCommon.Range("B156").Value = "Spare Parts (Solar):"

# ... then column C (Spanish translations) is filled in afterwards.
$wsCommon.Range("C151").Value = "Piezas De Repuesto (Eléctrico):"
$wsCommon.Range("C152").Value = "Piezas De Repuesto (Hardware):"
$wsCommon.Range("C153").Value = "Piezas De Repuesto (Vigilancia):"
$wsCommon.Range("C154").Value = "Piezas De Repuesto (Energía):"
$wsCommon.Range("C155").Value = "Piezas De Repuesto (Refrigeración):"
$wsCommon.Range("C156").Value = "Piezas De Repuesto (Solar):"

# Scroll the view down to the newly-added rows and select the cell right
# after the new data, then make this the active sheet/tab.
$excel.ActiveWindow.ScrollRow = 125
$wsCommon.Activate()
$wsCommon.Range("C157").Select()

# --- Workbook-level calculation settings ---
$excel.Iteration = $true
$excel.MaxChange = 0.0001
try { $excel.MultiThreadedCalculation.Enabled = $false } catch {}
